# Apply the "updated 4.0 files and mdl" revision:
#  - About!C1 date bumped to the new model-run date
#  - FPIEBP!B3:D3 (hard coal) balancing priorities re-ordered
#  - leave the cursor resting on FPIEBP!E3, matching the saved selection

$wb = $excel.ActiveWorkbook

$aboutSheet  = $wb.Worksheets.Item("About")
$fpiebpSheet = $wb.Worksheets.Item("FPIEBP")

# --- About sheet: bump the "last updated" date stamp (serial date value) ---
$aboutSheet.Range("C1").Value = 45379

# --- FPIEBP sheet: re-prioritize hard coal's production/imports/exports ---
$fpiebpSheet.Range("B3").Value = 1
$fpiebpSheet.Range("C3").Value = 3
$fpiebpSheet.Range("D3").Value = 2

# --- restore the active selection on the FPIEBP tab ---
$fpiebpSheet.Activate()
$fpiebpSheet.Range("E3").Select()
